$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "77.249.38"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "3.139.60"
$ws.Range("E3").Value = "  +5.48%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.03"
$ws.Range("E5").Value = "  +2.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "628.33"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  +13.09%  "

$ws.Range("E9").Value = "  +4.01%  "

$ws.Range("D10").Value = "3.142.77"
$ws.Range("E10").Value = "  +5.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.518"
$ws.Range("E11").Value = "  +19.36%  "

$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.38"
$ws.Range("E13").Value = "  +8.63%  "

$ws.Range("D14").Value = "3.719.06"
$ws.Range("E14").Value = "  +5.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000223"
$ws.Range("E15").Value = "  +19.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.53"
$ws.Range("E16").Value = "  +5.97%  "

$ws.Range("D17").Value = "77.119.14"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "3.141.96"
$ws.Range("E18").Value = "  +5.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.95"
$ws.Range("E19").Value = "  +3.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.21"
$ws.Range("E20").Value = "  +5.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.84"
$ws.Range("E21").Value = "  +26.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "412.49"
$ws.Range("E22").Value = "  +10.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.69"
$ws.Range("E23").Value = "  +9.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.66"
$ws.Range("E24").Value = "  +3.85%  "

$ws.Range("D25").Value = "3.304.49"
$ws.Range("E25").Value = "  +6.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.62"
$ws.Range("E26").Value = "  +8.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "74.84"
$ws.Range("E27").Value = "  +2.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.47"
$ws.Range("E28").Value = "  +8.59%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +7.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.994"
$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.72"
$ws.Range("E32").Value = "  +4.95%  "

$ws.Range("E33").Value = "  +5.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "520.39"
$ws.Range("E34").Value = "  +2.31%  "

$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.135"
$ws.Range("E36").Value = "  +20.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "21.71"
$ws.Range("E37").Value = "  +7.09%  "

$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "163.80"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.391"
$ws.Range("E40").Value = "  +1.66%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "195.53"
$ws.Range("E41").Value = "  +6.87%  "

$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.07"
$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.105"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.32"
$ws.Range("E45").Value = "  +8.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.796"
$ws.Range("E46").Value = "  +18.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.74"
$ws.Range("E47").Value = "  +6.42%  "

$ws.Range("E48").Value = "  +5.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "42.43"
$ws.Range("E49").Value = "  +0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.52"
$ws.Range("E50").Value = "  +10.30%  "

$ws.Range("E51").Value = "  +4.78%  "
